$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: add F3 = "Exp 7.png"
$ws.Range("F3").Value = "Exp 7.png"

# Row 4: fill in A4, B4, C4, F4 for Exp 8 data
$ws.Range("A4").Value = "Exp 8"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 10
$ws.Range("F4").Value = "Exp 8.png"

# Update selection to D6
$ws.Range("D6").Select()
